$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New build tag produced by the version-triggered build (fe11f8d1ec @ 2020-09-10 13:02),
# already present in the sheet at A4. Apply it to the remaining "test" rows that are
# getting a real smoke-test build number: ComparePages (A5), ModelPagesAccessories (A12),
# ModelPagesTechnology (A23) and SponsoredAthletes (A24).
$newBuild = "fe11f8d1ec built at 2020-09-10 13:02`n"

$targetRows = @(5, 12, 23, 24)
foreach ($r in $targetRows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $newBuild
    # Keep the row's default height - typing a multi-line value would otherwise
    # auto-grow the row, which isn't part of this edit.
    $ws.Rows($r).AutoFit()
}

# Matches the cursor coming to rest on the row right after the last edited cell (A24),
# as it would after typing the value into A24 and pressing Enter.
$ws.Range("A25").Select()
